# Swap the values of row 2 and row 4 for columns D, M, O, P, S
# (D=Fecha, M=Volumen, O=Precio máximo, P=Precio promedio ponderado, S=Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "O", "P", "S")

foreach ($col in $cols) {
    $r2 = $ws.Range($col + "2")
    $r4 = $ws.Range($col + "4")

    $v2 = $r2.Value2
    $v4 = $r4.Value2

    $r2.Value2 = $v4
    $r4.Value2 = $v2
}
